# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
#
# The "Date" column on the sheet holds a mislabeled value "7-1-2012-13" for
# every team row; it should read "2013-07-01" instead. Find the "Date"
# column header, then walk every data row underneath it and correct the
# label. The new value is entered with a leading apostrophe so Excel keeps
# it as literal text instead of re-parsing the ISO-looking string as a
# serial date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldLabel = "7-1-2012-13"
$newLabel = "2013-07-01"

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1
$firstCol  = $usedRange.Column
$lastCol   = $firstCol + $usedRange.Columns.Count - 1

$headerRow  = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($firstRow, $lastCol))
$dateHeader = $headerRow.Find("Date")
$dateCol    = $dateHeader.Column

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    if ($cell.Text -eq $oldLabel) {
        $cell.Value = "'" + $newLabel
    }
}
